# trafo_id -> gridnode_id refactor
#
# The "buildings" sheet has a header row (row 1) with one column per
# field. The column that used to be labelled "trafo_id" (column W) is
# renamed to "gridnode_id". Excel will automatically drop the now-unused
# "trafo_id" shared-string entry and append a new "gridnode_id" entry,
# which is exactly what happens in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the trafo_id header to gridnode_id (column W, row 1).
$ws.Range("W1").Value = "gridnode_id"

# Reflect the author's updated view/selection state: the sheet was
# scrolled so column Q is the left-most visible column, and the active
# selection moved from B5 to X6.
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 1
$ws.Range("X6").Select()

Write-Output "Renamed trafo_id -> gridnode_id and updated sheet selection"
